$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Address")

# Insert a new column before H (old H -> I, old I -> J) for the new
# "full state" lookup data ("shipping address and orders are done").
$ws.Columns.Item(8).Insert()

# Header
$ws.Range("H1").Value = "fulll_state"

# Body: full state name derived from the existing 2-letter state column (G)
$ws.Range("H2").Value = "Indiana"
$ws.Range("H3").Value = "Massachusetts"
$ws.Range("H4").Value = "Indiana"
$ws.Range("H5").Value = "Massachusetts"
$ws.Range("H6").Value = "New York"

# Match the column width used for the new column in the edited workbook
# (as close as this host's character-width quantization allows).
$ws.Columns.Item(8).ColumnWidth = 15.8

# Move the selection the way the author left it after the edit
$null = $ws.Range("I3").Select()
